# Weekly update: a new price record was reported for the week, which pushes
# the existing "Bruselas (repollito)" / "Vega Modelo de Temuco" series down
# by one row (the oldest row is re-appended at the bottom as row 130).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85; this shifts rows 85-129 down to 86-130
# (Excel also extends the used range / dimension to A1:R130 automatically).
$ws.Rows(85).Insert()

# Populate the newly inserted row 85 with this week's record.
$ws.Cells.Item(85, 1).Value = 10
$ws.Cells.Item(85, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value = "La Araucanía"
$ws.Cells.Item(85, 4).Value = 44813
$ws.Cells.Item(85, 5).Value = 9
$ws.Cells.Item(85, 6).Value = 100112035
$ws.Cells.Item(85, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 50
$ws.Cells.Item(85, 11).Value = 24000
$ws.Cells.Item(85, 12).Value = 25000
$ws.Cells.Item(85, 13).Value = 24400
$ws.Cells.Item(85, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(85, 15).Value = "Región Metropolitana"
$ws.Cells.Item(85, 16).Value = 2440
$ws.Cells.Item(85, 17).Value = 10
$ws.Cells.Item(85, 18).Value = "Hortaliza"
